# Adjusted excel-files for import tests.
#
# Replaces the placeholder "Ehto" condition-code values in column K (and
# clears/rewrites the now-unused localized "Ehto suomi/ruotsi/englanti"
# values in columns L:N) on rows 2-6 of Sheet1 with the real condition
# codes used by the automatic-fill import tests, and refreshes the active
# cell selection to reflect where the editor last left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: generic "muu" (other) condition code; localized texts unchanged.
$ws.Range("K2").Value = "muu Muu"

# Row 3: "lvm" (tuition fee) condition code; clears the old localized texts.
$ws.Range("K3").Value = "lvm Ehdollinen: lukuvuosimaksu maksettava määräaikaan mennessä, ennen kuin voit ilmoittautua"
$ws.Range("L3:N3").ClearContents()

# Row 4: "ttk" (degree-certificate copy) condition code; clears old texts.
$ws.Range("K4").Value = "ttk Ehdollinen: tutkintotodistuskopio hakuperusteena olleesta tutkinnosta toimitettava määräaikaan mennessä"
$ws.Range("L4:N4").ClearContents()

# Row 5: "ltt" (final degree certificate) condition code; clears old texts.
$ws.Range("K5").Value = "ltt Ehdollinen: lopullinen tutkintotodistus toimitettava määräaikaan mennessä"
$ws.Range("L5:N5").ClearContents()

# Row 6: generic "muu" (other) condition code with new "Testi" localized texts.
$ws.Range("K6").Value = "muu Muu"
$ws.Range("L6").Value = "Testi suomi"
$ws.Range("M6").Value = "Testi ruotsi"
$ws.Range("N6").Value = "Testi englanti"

# Move the active selection to where the editor left it.
$ws.Range("N7").Select()
